$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Text" number format on D-column cells whose new values would otherwise
# be auto-parsed by Excel as numbers (e.g. "73.94"), so they stay as plain text
# exactly like the rest of the Price column.
$textCells = @("D5","D6","D7","D9","D10","D12","D14","D15","D19","D21","D22","D23","D24","D27","D28","D30","D31","D32","D33","D36","D37","D39","D41","D42","D43","D45","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values row by row.
# Row 2
$ws.Range("D2").Value = '45.600.89'
$ws.Range("E2").Value = '  +6.06%  '
# Row 3
$ws.Range("D3").Value = '2.382.01'
# Row 4
$ws.Range("E4").Value = '  -0.83%  '
# Row 5
$ws.Range("D5").Value = '111.71'
$ws.Range("E5").Value = '  +6.25%  '
# Row 6
$ws.Range("D6").Value = '314.92'
$ws.Range("E6").Value = '  +1.83%  '
# Row 7
$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  +0.74%  '
# Row 8
$ws.Range("E8").Value = '  -0.10%  '
# Row 9
$ws.Range("D9").Value = '0.618'
$ws.Range("E9").Value = '  +2.18%  '
# Row 10
$ws.Range("D10").Value = '41.20'
$ws.Range("E10").Value = '  +3.67%  '
# Row 11
$ws.Range("E11").Value = '  +1.72%  '
# Row 12
$ws.Range("D12").Value = '8.57'
$ws.Range("E12").Value = '  +3.48%  '
# Row 13
$ws.Range("E13").Value = '  +2.03%  '
# Row 14
$ws.Range("D14").Value = '0.989'
$ws.Range("E14").Value = '  +0.30%  '
# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '15.60'
$ws.Range("E15").Value = '  +2.03%  '
# Row 16
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.745.62'
$ws.Range("E16").Value = '  +3.52%  '
# Row 17
$ws.Range("D17").Value = '2.379.94'
$ws.Range("E17").Value = '  +3.25%  '
# Row 18
$ws.Range("D18").Value = '45.511.17'
$ws.Range("E18").Value = '  +6.66%  '
# Row 19
$ws.Range("D19").Value = '7.35'
$ws.Range("E19").Value = '  +0.30%  '
# Row 20
$ws.Range("E20").Value = '  +1.83%  '
# Row 21
$ws.Range("D21").Value = '13.13'
$ws.Range("E21").Value = '  -5.03%  '
# Row 22
$ws.Range("D22").Value = '73.94'
# Row 23
$ws.Range("D23").Value = '3.49'
$ws.Range("E23").Value = '  +1.37%  '
# Row 24
$ws.Range("D24").Value = '262.70'
$ws.Range("E24").Value = '  -2.06%  '
# Row 25
$ws.Range("E25").Value = '  +3.22%  '
# Row 26
$ws.Range("E26").Value = '  -0.43%  '
# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '11.16'
$ws.Range("E27").Value = '  +2.04%  '
# Row 28
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = '7.48'
$ws.Range("E28").Value = '  -2.90%  '
# Row 29
$ws.Range("E29").Value = '  +2.47%  '
# Row 30
$ws.Range("D30").Value = '22.59'
$ws.Range("E30").Value = '  +1.64%  '
# Row 31
$ws.Range("D31").Value = '38.22'
$ws.Range("E31").Value = '  +1.17%  '
# Row 32
$ws.Range("D32").Value = '0.0968'
$ws.Range("E32").Value = '  +11.81%  '
# Row 33
$ws.Range("D33").Value = '170.72'
$ws.Range("E33").Value = '  +3.10%  '
# Row 34
$ws.Range("E34").Value = '  +4.28%  '
# Row 35
$ws.Range("E35").Value = '  +0.36%  '
# Row 36
$ws.Range("D36").Value = '0.117'
$ws.Range("E36").Value = '  +3.68%  '
# Row 37
$ws.Range("D37").Value = '4.83'
$ws.Range("E37").Value = '  +4.22%  '
# Row 38
$ws.Range("E38").Value = '  +10.97%  '
# Row 39
$ws.Range("D39").Value = '2.99'
$ws.Range("E39").Value = '  +7.28%  '
# Row 40
$ws.Range("E40").Value = '  +0.28%  '
# Row 41
$ws.Range("D41").Value = '1.74'
$ws.Range("E41").Value = '  +10.92%  '
# Row 42
$ws.Range("D42").Value = '102.91'
$ws.Range("E42").Value = '  -4.95%  '
# Row 43
$ws.Range("D43").Value = '0.236'
$ws.Range("E43").Value = '  +3.61%  '
# Row 44
$ws.Range("E44").Value = '  +8.62%  '
# Row 45
$ws.Range("D45").Value = '70.36'
$ws.Range("E45").Value = '  -1.44%  '
# Row 46
$ws.Range("E46").Value = '  +0.15%  '
# Row 47
$ws.Range("D47").Value = '84.28'
$ws.Range("E47").Value = '  +10.87%  '
# Row 48
$ws.Range("D48").Value = '113.79'
$ws.Range("E48").Value = '  +1.81%  '
# Row 49
$ws.Range("D49").Value = '9.38'
$ws.Range("E49").Value = '  +6.00%  '
# Row 50
$ws.Range("D50").Value = '5.56'
$ws.Range("E50").Value = '  +7.37%  '
# Row 51
$ws.Range("D51").Value = '1.649.62'
$ws.Range("E51").Value = '  -2.86%  '
